# "run prepare & render with final data"
# The underlying data-prep pipeline was re-run: the four "Europe *" poll-breakdown
# columns were dropped entirely, two rows of statements were reordered, and every
# percentage in the table was refreshed with the latest render output.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the 4 "Europe ..." breakdown columns (D:G); remaining columns shift left.
$ws.Range("D1:G11").EntireColumn.Delete() | Out-Null

# Header row (A1 is already an empty string and stays untouched).
$ws.Range("B1").Value = "`$ bold('All')"
$ws.Range("C1").Value = "Millionaires"
$ws.Range("D1").Value = "Japan Non-voters"
$ws.Range("E1").Value = "Japan Left"
$ws.Range("F1").Value = "Japan Center/Right"
$ws.Range("G1").Value = "Saudi Arabia"
$ws.Range("H1").Value = "Saudi citizens"
$ws.Range("I1").Value = "U.S. Non-voters"
$ws.Range("J1").Value = "U.S. Harris"
$ws.Range("K1").Value = "U.S. Trump"

# Re-populate rows 2-11 with the statement labels (reordered) and refreshed data.
# Row 2
$ws.Range("A2").Value = "Minimum tax of 2% on billionaires'`nwealth, in voluntary countries"
$ws.Range("B2").Value = 0.633130423220635
$ws.Range("C2").Value = 0.556791704304188
$ws.Range("D2").Value = 0.493056521315741
$ws.Range("E2").Value = 0.608803357414973
$ws.Range("F2").Value = 0.531182230184829
$ws.Range("G2").Value = 0.670833846908728
$ws.Range("H2").Value = 0.692812409550454
$ws.Range("I2").Value = 0.549295978763342
$ws.Range("J2").Value = 0.807022646248167
$ws.Range("K2").Value = 0.452892965198557
$ws.Rows.Item(2).AutoFit() | Out-Null
# Row 3
$ws.Range("A3").Value = "Bridgetown initiative: MDBs expanding sustainable`ninvestments in LICs, and at lower interest rates"
$ws.Range("B3").Value = 0.563498853260559
$ws.Range("C3").Value = 0.592539136795488
$ws.Range("D3").Value = 0.383411391668753
$ws.Range("E3").Value = 0.510153602277153
$ws.Range("F3").Value = 0.465702404553403
$ws.Range("G3").Value = 0.699951386105416
$ws.Range("H3").Value = 0.705866678710194
$ws.Range("I3").Value = 0.449570088249459
$ws.Range("J3").Value = 0.735968083342179
$ws.Range("K3").Value = 0.386344486755095
$ws.Rows.Item(3).AutoFit() | Out-Null
# Row 4
$ws.Range("A4").Value = "L&D: Developed countries financing a fund to help`nvulnerable countries cope with climate Loss and damage"
$ws.Range("B4").Value = 0.548888702456242
$ws.Range("C4").Value = 0.528919523634788
$ws.Range("D4").Value = 0.336693874193082
$ws.Range("E4").Value = 0.505578632976523
$ws.Range("F4").Value = 0.481053993918987
$ws.Range("G4").Value = 0.754595050954226
$ws.Range("H4").Value = 0.764918944465133
$ws.Range("I4").Value = 0.455662843526806
$ws.Range("J4").Value = 0.697789663211793
$ws.Range("K4").Value = 0.345516117316969
$ws.Rows.Item(4).AutoFit() | Out-Null
# Row 5
$ws.Range("A5").Value = "Debt relief for vulnerable countries, suspending`npayments until they are more able to repay"
$ws.Range("B5").Value = 0.492515087699993
$ws.Range("C5").Value = 0.427704966662765
$ws.Range("D5").Value = 0.311581900605115
$ws.Range("E5").Value = 0.460084045480565
$ws.Range("F5").Value = 0.38441921408349
$ws.Range("G5").Value = 0.703571225719179
$ws.Range("H5").Value = 0.7538556364906
$ws.Range("I5").Value = 0.446547267517491
$ws.Range("J5").Value = 0.597618388012403
$ws.Range("K5").Value = 0.352144673086557
$ws.Rows.Item(5).AutoFit() | Out-Null
# Row 6
$ws.Range("A6").Value = "At least 0.7% of developed countries' GDP in foreign aid"
$ws.Range("B6").Value = 0.487666417243067
$ws.Range("C6").Value = 0.507683351069014
$ws.Range("D6").Value = 0.219184379605566
$ws.Range("E6").Value = 0.387322696004644
$ws.Range("F6").Value = 0.377769238751658
$ws.Range("G6").Value = 0.688471530144484
$ws.Range("H6").Value = 0.721961860164294
$ws.Range("I6").Value = 0.361562972247251
$ws.Range("J6").Value = 0.632667573918023
$ws.Range("K6").Value = 0.339453064229604
$ws.Rows.Item(6).AutoFit() | Out-Null
# Row 7
$ws.Range("A7").Value = "Raise global minimum tax on profit from 15% to 35%,`nallocating revenues to countries based on sales"
$ws.Range("B7").Value = 0.486367361421124
$ws.Range("C7").Value = 0.518104168638399
$ws.Range("D7").Value = 0.31357267707136
$ws.Range("E7").Value = 0.484806208512163
$ws.Range("F7").Value = 0.448145256545193
$ws.Range("G7").Value = 0.530644213623694
$ws.Range("H7").Value = 0.591894601389567
$ws.Range("I7").Value = 0.366262634781694
$ws.Range("J7").Value = 0.635651938591609
$ws.Range("K7").Value = 0.33853677151088
$ws.Rows.Item(7).AutoFit() | Out-Null
# Row 8
$ws.Range("A8").Value = "NCQG: Developing countries providing `$300 bn a`nyear in climate finance for developing countries"
$ws.Range("B8").Value = 0.484425269846573
$ws.Range("C8").Value = 0.452900311678667
$ws.Range("D8").Value = 0.197887265448319
$ws.Range("E8").Value = 0.385404090605957
$ws.Range("F8").Value = 0.379724102910581
$ws.Range("G8").Value = 0.672464791241274
$ws.Range("H8").Value = 0.71838912281039
$ws.Range("I8").Value = 0.378924737497152
$ws.Range("J8").Value = 0.619589144561117
$ws.Range("K8").Value = 0.258254482209228
$ws.Rows.Item(8).AutoFit() | Out-Null
# Row 9
$ws.Range("A9").Value = "International levy on shipping carbon emissions,`nreturned to countries based on population"
$ws.Range("B9").Value = 0.472274429131213
$ws.Range("C9").Value = 0.476364502493748
$ws.Range("D9").Value = 0.223394314273261
$ws.Range("E9").Value = 0.326906105388108
$ws.Range("F9").Value = 0.340508401165705
$ws.Range("G9").Value = 0.603725863836496
$ws.Range("H9").Value = 0.661498159140946
$ws.Range("I9").Value = 0.369259741689718
$ws.Range("J9").Value = 0.605720762515847
$ws.Range("K9").Value = 0.349545584826317
$ws.Rows.Item(9).AutoFit() | Out-Null
# Row 10
$ws.Range("A10").Value = "Expand Security Council to new permanent members (e.g.`nIndia, Brazil, African Union), restrict veto use"
$ws.Range("B10").Value = 0.463687229299451
$ws.Range("C10").Value = 0.497513940092693
$ws.Range("D10").Value = 0.231140663914579
$ws.Range("E10").Value = 0.39288890608071
$ws.Range("F10").Value = 0.400497772491134
$ws.Range("G10").Value = 0.629350439518224
$ws.Range("H10").Value = 0.677149008732137
$ws.Range("I10").Value = 0.355170865465525
$ws.Range("J10").Value = 0.611945859008054
$ws.Range("K10").Value = 0.310545566749026
$ws.Rows.Item(10).AutoFit() | Out-Null
# Row 11
$ws.Range("A11").Value = "International levy on aviation carbon emissions, raising`nprices by 30%, returned to countries based on population"
$ws.Range("B11").Value = 0.373752935747861
$ws.Range("C11").Value = 0.367935781075227
$ws.Range("D11").Value = 0.206480752740589
$ws.Range("E11").Value = 0.301394636393387
$ws.Range("F11").Value = 0.28022832538204
$ws.Range("G11").Value = 0.533333498726061
$ws.Range("H11").Value = 0.58586931877385
$ws.Range("I11").Value = 0.285288177307175
$ws.Range("J11").Value = 0.48529662511572
$ws.Range("K11").Value = 0.257067972351958
$ws.Rows.Item(11).AutoFit() | Out-Null
